$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.848.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.81%  '
$ws.Range("D3").Value = "'2.416.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.69%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'554.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.69%  '
$ws.Range("D6").Value = "'138.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.99%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("E9").Value = '  +4.91%  '
$ws.Range("D10").Value = "'5.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.39%  '
$ws.Range("E11").Value = '  +1.65%  '
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = "'24.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'2.848.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").Value = "'59.742.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("E16").Value = '  +4.38%  '
$ws.Range("D17").Value = "'2.438.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.52%  '
$ws.Range("E18").Value = '  +6.43%  '
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("D20").Value = "'333.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = "'6.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = "'64.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.88%  '
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").Value = "'0.0₃0786"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.82%  '
$ws.Range("E29").Value = '  +3.37%  '
$ws.Range("D30").Value = "'170.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("D32").Value = "'18.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("E35").Value = '  +5.26%  '
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = "'1.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("D39").Value = "'40.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("D40").Value = "'0.423"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.96%  '
$ws.Range("D41").Value = "'312.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.42%  '
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("D43").Value = "'143.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").Value = "'0.0963"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("E45").Value = '  +4.30%  '
$ws.Range("D46").Value = "'0.418"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.35%  '
$ws.Range("D47").Value = "'19.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("E49").Value = '  +2.94%  '
$ws.Range("D50").Value = "'11.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("E51").Value = '  +4.73%  '
